# [IMPL] CRUD ADMINISTRADOR (ADMIN)
# Adds the 5 new "administrador" CRUD tasks into the SPRINT-BACKLOG sheet
# (sprint 3 section), right above the existing "Introducción ..." rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT-BACKLOG")

# Insert 5 fresh rows right before the current row 77 - this pushes the
# existing rows 77-247 down to 82-252 (formulas on the AUX sheet that
# reference 'SPRINT-BACKLOG'!...40:...81 auto-expand to ...40:...86, and the
# blank "filler" rows at the bottom of the used range grow from 247 to 252).
$ws.Rows("77:81").Insert()

# Copy the formatting (styles/number formats/borders) of the row right below
# (the one that used to be row 77, now row 82) onto the 5 new blank rows, so
# they reuse the same style ids as every other data row instead of picking up
# a freshly interpolated style.
$ws.Range("A82:G82").Copy()
$ws.Range("A77:G81").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new description strings are entered in this specific order (matching
# the order they land in xl/sharedStrings.xml) so the shared-string table
# indices line up with the source edit: mostrar, editar, insertar, listar,
# eliminar - not the top-to-bottom row order.
$ws.Range("A78").Value = "Implementación de mostrar los detalles de un administrador (ADMIN)"
$ws.Range("A80").Value = "Implementación de la funcionalidad de editar los detalles de un administrador(ADMIN) "
$ws.Range("A79").Value = "Implementación de la funcionalidad de insertar un nuevo administrador (ADMIN) "
$ws.Range("A77").Value = "AImplementación de la funcionalidad de listar administradores (ADMIN)"
$ws.Range("A81").Value = "Implementación de la funcionalidad de eliminar un administrador (ADMIN)"

# Row 77: listar administradores
$ws.Range("B77").Value = 1
$ws.Range("C77").Value = ""
$ws.Range("D77").Value = 44281
$ws.Range("E77").Value = ""
$ws.Range("F77").Value = "REALIZADO"

# Row 78: mostrar los detalles de un administrador
$ws.Range("B78").Value = 1
$ws.Range("C78").Value = 44274
$ws.Range("D78").Value = 44281
$ws.Range("E78").Value = 44277
$ws.Range("F78").Value = "REALIZADO"

# Row 79: insertar un nuevo administrador
$ws.Range("B79").Value = 1
$ws.Range("C79").Value = ""
$ws.Range("D79").Value = 44281
$ws.Range("E79").Value = ""
$ws.Range("F79").Value = "REALIZADO"

# Row 80: editar los detalles de un administrador
$ws.Range("B80").Value = 1
$ws.Range("C80").Value = ""
$ws.Range("D80").Value = 44281
$ws.Range("E80").Value = ""
$ws.Range("F80").Value = "REALIZADO"

# Row 81: eliminar un administrador
$ws.Range("B81").Value = 1
$ws.Range("C81").Value = 44274
$ws.Range("D81").Value = 44281
$ws.Range("E81").Value = 44277
$ws.Range("F81").Value = "REALIZADO"

# Let Excel auto-fit the new rows' heights from their wrapped text, same as
# every other row in this sheet (none of them carry an explicit
# customHeight="1").
$ws.Rows("77:81").EntireRow.AutoFit()

$wb.Save()
